# Update indicator 6.3.1 metadata sheet:
#  - refresh the indicator name (now also covers industrial wastewater)
#  - refresh the responsible organization / contact person / phone / website
#    (the sheet's reporting department and contact person changed)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Информация об индикаторе -> Индикатор
$ws.Range("B4").Value2 = "6.3.1. Доля безопасно очищаемых хозяйственнобытовых и промышленных сточных вод "
$ws.Range("B4").WrapText = $true

# 2. Информация об организации -> Организация
$ws.Range("B6").Value2 = "Национальный статистический комитет Кыргызской Республики`nУправление цифрового развития и статистики устойчивого развития"

# 2. Информация об организации -> Контактное лицо (лица) / Координатор
$ws.Range("B7").Value2 = "Мамбеталиев Т.А."

# 2. Информация об организации -> Электронная почта контактного лица
$ws.Range("B8").Value2 = "Sdg_nsc@stat.kg "

# 2. Информация об организации -> Телефон контактного лица
$ws.Range("B9").Value2 = "(0312) 62 56 07"

# 2. Информация об организации -> Сайт организации (если есть)
$ws.Range("B10").Value2 = "www.stat.gov.kg"

# Leave the selection where the author's last edit landed
[void]$ws.Range("B9").Select()
